$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '29.868.24'
$ws.Range('E2').Value = '  -1.61%  '
$ws.Range('D3').Value = '1.887.62'
$ws.Range('E3').Value = '  -2.51%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.19%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '0.7244'
$ws.Range('E5').Value = '  -6.42%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '241.76'
$ws.Range('E6').Value = '  -1.84%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.003'
$ws.Range('E7').Value = '  +0.20%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.3078'
$ws.Range('E8').Value = '  -3.85%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '26.06'
$ws.Range('E9').Value = '  -6.19%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.06855'
$ws.Range('E10').Value = '  -2.85%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.07945'
$ws.Range('E11').Value = '  -0.96%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '0.7619'
$ws.Range('E12').Value = '  -2.29%  '
$ws.Range('D13').Value = '1.871.86'
$ws.Range('E13').Value = '  -3.27%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '5.221'
$ws.Range('E14').Value = '  -2.57%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '90.80'
$ws.Range('E15').Value = '  -4.09%  '
$ws.Range('D16').Value = '29.891.79'
$ws.Range('E16').Value = '  -1.51%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '14.01'
$ws.Range('E17').Value = '  -3.07%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '5.728'
$ws.Range('E18').Value = '  -1.20%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '238.82'
$ws.Range('E19').Value = '  -6.59%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.000007715'
$ws.Range('E20').Value = '  -3.01%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '1.004'
$ws.Range('E21').Value = '  +0.26%  '
$ws.Range('D22').Value = '2.106.54'
$ws.Range('E22').Value = '  -3.77%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '1.006'
$ws.Range('E23').Value = '  +0.29%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '6.787'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '9.238'
$ws.Range('E25').Value = '  -3.29%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '164.61'
$ws.Range('E26').Value = '  +0.56%  '
$ws.Range('E27').Value = '  -0.85%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.1270'
$ws.Range('E28').Value = '  -5.97%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '2.001'
$ws.Range('E29').Value = '  -11.98%  '
$ws.Range('E30').Value = '  -0.93%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.527'
$ws.Range('E31').Value = '  +0.76%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '4.280'
$ws.Range('E32').Value = '  -3.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '4.041'
$ws.Range('E33').Value = '  -1.78%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '0.05043'
$ws.Range('E34').Value = '  -2.19%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.262'
$ws.Range('E35').Value = '  -1.58%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.7287'
$ws.Range('E36').Value = '  -2.47%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '2.729'
$ws.Range('E37').Value = '  -1.87%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.01906'
$ws.Range('E38').Value = '  -2.31%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '2.765'
$ws.Range('E39').Value = '  -1.57%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '6.300'
$ws.Range('E40').Value = '  -2.01%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '74.05'
$ws.Range('E41').Value = '  -5.75%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '0.4402'
$ws.Range('E42').Value = '  -2.20%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.910'
$ws.Range('E43').Value = '  -3.06%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '1.004'
$ws.Range('E44').Value = '  +0.15%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.8343'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '100.56'
$ws.Range('E46').Value = '  -0.21%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '7.574'
$ws.Range('E47').Value = '  +0.84%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '9.731'
$ws.Range('E48').Value = '  -0.87%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '37.26'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').Value = '2.028.96'
$ws.Range('E50').Value = '  -2.68%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '936.35'
$ws.Range('E51').Value = '  -5.03%  '
